$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 "blog" slot cells get new ser numbers (C9 gets the brand-new ser:141 post;
# H9/J9/L9 shift to the next ser values as the feed advances).
$ws.Range("C9").Value = "type: blog`nwidth: 2`nheight: 1`nser: 141"
$ws.Range("H9").Value = "type: blog`nwidth: 2`nheight: 1`nser: 139"
$ws.Range("J9").Value = "type: blog`nwidth: 2`nheight: 1`nser: 137"
$ws.Range("L9").Value = "type: blog`nwidth: 2`nheight: 1`nser: 136"

# Update the active view/selection to C9.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$win.ScrollRow = 9
$null = $ws.Range("C9").Select()
